$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column D (us_businessvalue) values per row as per the recorded diff.
$values = @{
    2 = 89
    3 = 89
    4 = 89
    5 = 89
    6 = 89
    7 = 89
    8 = 89
    9 = 89
    10 = 5
    11 = 21
    12 = 5
    13 = 5
    14 = 21
    15 = 34
    16 = 5
    17 = 34
    18 = 5
    19 = 5
    20 = 5
    21 = 5
    22 = 5
    23 = 5
    24 = 5
    25 = 5
    26 = 5
    27 = 5
    28 = 13
    29 = 13
    30 = 13
    31 = 13
    32 = 13
    33 = 13
    34 = 13
    35 = 13
    36 = 13
    37 = 13
    38 = 13
    39 = 13
    40 = 13
    41 = 13
    42 = 13
    43 = 13
    44 = 13
    45 = 13
    46 = 13
    47 = 13
    48 = 13
    49 = 13
    50 = 13
    51 = 13
    52 = 13
    53 = 13
    54 = 13
    55 = 13
    56 = 13
    57 = 8
    58 = 8
    59 = 8
    60 = 5
    61 = 8
    62 = 34
    63 = 5
    64 = 8
    65 = 89
    66 = 8
    67 = 13
    68 = 21
    69 = 5
    70 = 13
    71 = 89
    72 = 8
    73 = 8
    74 = 21
    75 = 13
    76 = 55
    77 = 13
    78 = 89
    79 = 13
    80 = 5
    81 = 8
    82 = 8
    83 = 8
    84 = 21
    85 = 89
    87 = 55
    88 = 55
    89 = 55
    90 = 55
    91 = 55
    92 = 55
    93 = 55
    94 = 55
    95 = 55
    96 = 55
    97 = 55
    98 = 55
    99 = 13
    100 = 5
    101 = 13
    102 = 21
    103 = 8
    104 = 21
    105 = 89
    106 = 8
    107 = 8
    108 = 8
    109 = 5
    110 = 13
    111 = 5
    112 = 8
    113 = 8
    114 = 34
    115 = 34
    116 = 34
    117 = 34
    118 = 5
    119 = 8
    120 = 8
    121 = 13
    122 = 8
    123 = 8
    124 = 8
    125 = 13
    126 = 34
    127 = 34
    128 = 13
    129 = 13
    130 = 13
    131 = 13
    132 = 13
    133 = 13
    134 = 89
    135 = 5
    136 = 8
    137 = 13
    138 = 13
    139 = 5
    140 = 89
    141 = 5
    142 = 5
    143 = 5
    144 = 5
    145 = 13
    146 = 13
    147 = 13
    148 = 13
    149 = 13
    150 = 13
    151 = 89
    152 = 21
    153 = 8
    154 = 8
    155 = 8
    156 = 5
    157 = 13
    158 = 5
    159 = 5
    160 = 5
    161 = 5
    162 = 5
    163 = 5
    164 = 5
    165 = 5
    166 = 5
    167 = 5
    168 = 5
    169 = 5
    170 = 5
    171 = 5
    172 = 5
    173 = 5
    174 = 5
    175 = 5
    176 = 5
    177 = 5
    178 = 5
    179 = 34
    180 = 5
    181 = 8
    182 = 5
    183 = 5
    184 = 5
    185 = 5
    186 = 5
    187 = 5
    188 = 5
    189 = 5
    190 = 5
    191 = 5
    192 = 5
    193 = 5
    194 = 5
    195 = 5
    196 = 5
    197 = 5
    198 = 5
    199 = 5
    200 = 5
    201 = 5
    202 = 5
    203 = 5
    204 = 5
    205 = 5
    206 = 5
    207 = 5
    208 = 5
    209 = 5
    210 = 5
    211 = 89
    212 = 13
    213 = 21
    214 = 5
    215 = 5
    216 = 8
    217 = 55
    218 = 5
    219 = 8
    220 = 8
    221 = 5
    222 = 8
    223 = 8
    224 = 8
    225 = 8
    226 = 21
    227 = 21
    228 = 21
    229 = 5
    230 = 5
    231 = 5
    232 = 34
    233 = 21
    234 = 13
    235 = 13
    236 = 55
    237 = 13
    238 = 5
    239 = 34
    240 = 5
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}

